# Updated code and test data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "DQ_Report"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Update existing rows 2-11 (columns A,B,C) to their new values.
$ws1.Cells.Item(2,1).Value = "P_20085751"
$ws1.Cells.Item(2,2).Value = "G70"
$ws1.Cells.Item(2,3).Value = 586

$ws1.Cells.Item(3,1).Value = "P_20085752"
$ws1.Cells.Item(3,2).Value = "G70"
$ws1.Cells.Item(3,3).Value = 589

$ws1.Cells.Item(4,1).Value = "P_20085754"
$ws1.Cells.Item(4,2).Value = "E84.80"
$ws1.Cells.Item(4,3).Value = 588

$ws1.Cells.Item(5,1).Value = "P_20085755"
$ws1.Cells.Item(5,2).Value = "E75.2"
$ws1.Cells.Item(5,3).Value = 325

$ws1.Cells.Item(6,1).Value = "P_20085756"
$ws1.Cells.Item(6,2).Value = "E75.2"
$ws1.Cells.Item(6,3).Value = 320

$ws1.Cells.Item(7,1).Value = "P_20085757"
$ws1.Cells.Item(7,3).Value = 586

$ws1.Cells.Item(8,1).Value = "P_20085758"
$ws1.Cells.Item(8,3).Value = 587

$ws1.Cells.Item(9,1).Value = "P_20085759"
$ws1.Cells.Item(9,2).Value = "E84.0"
$ws1.Cells.Item(9,3).Value = ""

$ws1.Cells.Item(10,1).Value = "P_20085760"
$ws1.Cells.Item(10,2).Value = "D45"

$ws1.Cells.Item(11,1).Value = "P_20085761"
$ws1.Cells.Item(11,2).Value = ""

# Remove column D (dq_msg) entirely - header + all values.
$ws1.Columns.Item(4).Delete()

# Append the new rows (12-16).
$ws1.Cells.Item(12,1).Value = "P_20085762"
$ws1.Cells.Item(12,2).Value = "E66.89"
$ws1.Cells.Item(12,3).Value = 320

$ws1.Cells.Item(13,1).Value = "P_20085764"
$ws1.Cells.Item(13,2).Value = "E66.89"

$ws1.Cells.Item(14,1).Value = "P_20085764"
$ws1.Cells.Item(14,2).Value = "E75.2"

$ws1.Cells.Item(15,1).Value = "P_20085767"
$ws1.Cells.Item(15,2).Value = "E85.0"
$ws1.Cells.Item(15,3).Value = 586

$ws1.Cells.Item(16,1).Value = "P_20085770"
$ws1.Cells.Item(16,2).Value = "J09"

# ---------------------------------------------------------------
# Sheet 2: "Statistik"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2,2).Value = 15.47
$ws2.Cells.Item(2,3).Value = 84.53
$ws2.Cells.Item(2,4).Value = 88.23999999999999
$ws2.Cells.Item(2,5).Value = 97.40000000000001
$ws2.Cells.Item(2,6).Value = 34
$ws2.Cells.Item(2,7).Value = 305
$ws2.Cells.Item(2,8).Value = 395

$ws2.Cells.Item(1,9).Value = "case_no"
$ws2.Cells.Item(2,9).Value = 420

# ---------------------------------------------------------------
# Sheet 3: "Projectathon"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(3,6).Value = 1
$ws3.Cells.Item(3,7).Value = 25

$ws3.Cells.Item(6,1).Value = "C2"
$ws3.Cells.Item(6,2).Value = "Anderson‐Fabry‐Krankheit"
$ws3.Cells.Item(6,3).Value = 324
$ws3.Cells.Item(6,4).Value = "E75.2"
$ws3.Cells.Item(6,5).Value = 9
$ws3.Cells.Item(6,6).Value = 0
$ws3.Cells.Item(6,7).Value = 0
